# Insert a new pricing record for Membrillo (Vega Modelo de Temuco) at row 273.
# Excel's EntireRow.Insert shifts row 273 (and everything below it) down by
# one, which is exactly the shape described by the diff (rows 273-298 ->
# 274-299, each row keeping its original data) while the new row 273 holds
# the freshly reported week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(273).EntireRow.Insert()

$ws.Cells.Item(273, 1).Value = 10
$ws.Cells.Item(273, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(273, 3).Value = "La Araucanía"
$ws.Cells.Item(273, 4).Value = 45106
$ws.Cells.Item(273, 5).Value = 9
$ws.Cells.Item(273, 6).Value = "Fruta"
$ws.Cells.Item(273, 7).Value = 100104
$ws.Cells.Item(273, 8).Value = "Frutos de pepita"
$ws.Cells.Item(273, 9).Value = 100104003
$ws.Cells.Item(273, 10).Value = "Membrillo"
$ws.Cells.Item(273, 11).Value = "Champion"
$ws.Cells.Item(273, 12).Value = "Primera"
$ws.Cells.Item(273, 13).Value = 80
$ws.Cells.Item(273, 14).Value = 14000
$ws.Cells.Item(273, 15).Value = 14000
$ws.Cells.Item(273, 16).Value = 14000
$ws.Cells.Item(273, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(273, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(273, 19).Value = 778
$ws.Cells.Item(273, 20).Value = 18
